# Adding the changes we made on may 9th
#
# Inserts 11 new accelerometer samples right after the header row (so they
# become the new rows 2-12), pushing the previously-existing data down by
# 11 rows. The data set also drops its final (oldest) sample, which falls
# out of the new A1:C31 used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of x/y/z data to insert right below the header row.
$newRows = @(
    @(-2.507025241851806, 9.102962493896484, -0.6246470808982849),
    @(-2.812539100646973, 9.111927032470703, -0.3471660315990448),
    @(-2.887884616851806, 9.205938339233398, -0.2687076330184936),
    @(-2.862497329711914, 9.105484008789062, -0.2620421051979065),
    @(-2.736449241638184, 9.052282333374023, -0.2773746848106384),
    @(-2.651521682739258, 9.127286911010742, 0.0024068877100944),
    @(-2.57433032989502,  9.058347702026367, 0.0778785794973373),
    @(-2.764423370361328, 9.000140190124512, -0.2859586775302887),
    @(-2.818140029907227, 8.811227798461914, -0.6212533712387085),
    @(-3.191051483154297, 9.015185356140137, -0.4772885143756866),
    @(-2.787490367889404, 8.648155212402344, -1.196338057518005)
)

$insertCount = $newRows.Count

# Shift the existing data (rows 2..21) down by inserting $insertCount blank
# rows right after the header row (row 1).
$insertRange = $ws.Range("A2:C$($insertCount + 1)")
$insertRange.EntireRow.Insert()

# Row-insert copies formatting down from the row above (the bold header),
# so strip that back off - the data rows carry no explicit style.
$ws.Range("A2:C$($insertCount + 1)").ClearFormats()

# Write the new sample rows into the freshly inserted rows (2..12).
for ($i = 0; $i -lt $insertCount; $i++) {
    $r = 2 + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# The original last row (previously row 21, now shifted to row 32) is no
# longer part of the data set - remove it so the used range ends at row 31.
$lastRow = 21 + $insertCount
$ws.Range("A$($lastRow):C$($lastRow)").EntireRow.Delete()
